# PYME-4531: Add "Approved At" column (date the request was approved) to the
# european_funds_requests export template, right after "Created At".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D, pushing the existing D..M columns one place to the
# right (Customer ID, Technical Email, Last Change At, ... Domain).
$ws.Range("D1").EntireColumn.Insert()

# Match the width/grouping of the neighbouring "Created At" column (C) that
# the new column was inserted next to.
$ws.Columns.Item(4).ColumnWidth = 18.75
$ws.Columns.Item(4).OutlineLevel = 1

# New header label + shared string.
$ws.Range("D1").Value = "Approved At"

# The autofilter now needs to span through the new last column (N).
$ws.AutoFilterMode = $false
$ws.Range("A1:N1").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the autofilter.
$wb.Names.Item("Data!_FilterDatabase").RefersTo = "=Data!`$A`$1:`$N`$1"

# Put the active selection on the newly added header cell.
$ws.Range("D1").Select() | Out-Null
